# Fixed s bug in stats2
# The rows of the stats table (A3:F23) had their values shuffled into the
# wrong rows. This restores each row's data (A..F) to the correct row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(501, 9, 52, 30, 75, 45),    # row 3
    @(201, 9, 30, 15, 45, 30),    # row 4
    @(1001, 18, 30, 75, 60, 72),  # row 5
    @(301, 6, 45, 30, 60, 45),    # row 6
    @(701, 3, 90, 45, 97, 15),    # row 7
    @(601, 9, 60, 67, 60, 42),    # row 8
    @(1201, 2, 10, 10, 10, 10),   # row 9
    @(1203, 3, 15, 15, 15, 15),   # row 10
    @(101, 9, 30, 15, 60, 15),    # row 11
    @(1202, 2, 10, 10, 10, 10),   # row 12
    @(902, 1, 0, 0, 0, 0),        # row 13
    @(401, 9, 48, 67, 75, 45),    # row 14
    @(801, 3, 67, 65, 52, 45),    # row 15
    @(502, 0, 4, 0, 0, 0),        # row 16
    @(1, 0, 2, 2, 2, 2),          # row 17
    @(2, 0, 2, 2, 2, 2),          # row 18
    @(3, 0, 3, 3, 3, 3),          # row 19 (unchanged)
    @(802, 0, 4, 5, 4, 0),        # row 20
    @(1101, 0, 15, 30, 30, 0),    # row 21
    @(602, 0, 0, 4, 0, 9),        # row 22
    @(402, 0, 0, 4, 0, 0)         # row 23
)

$startRow = 3
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowVals[$c]
    }
}
